$wb = $excel.ActiveWorkbook

# Rename the 'General' sheet to 'Table'
$ws = $wb.Worksheets.Item("General")
$ws.Name = "Table"

# Update the selection on that sheet to R85
$ws.Activate()
$ws.Range("R85").Select()
